$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value2 = 13.70755662053669
$ws.Range("D2").Value2 = 8.387745344476244
$ws.Range("E2").Value2 = 14.33178780110785
$ws.Range("F2").Value2 = 40.1508528792226
$ws.Range("G2").Value2 = 3.706096110877508
$ws.Range("J2").Value2 = 11.02464558371742
$ws.Range("K2").Value2 = 20.20678921467036
$ws.Range("L2").Value2 = 9.739148067921626
$ws.Range("O2").Value2 = 30.84979118004521
$ws.Range("C3").Value2 = 13.66056435596128
$ws.Range("D3").Value2 = 8.337492593953266
$ws.Range("E3").Value2 = 14.32582874220938
$ws.Range("F3").Value2 = 40.33034539568528
$ws.Range("G3").Value2 = 3.70868512356346
$ws.Range("J3").Value2 = 11.05092804880945
$ws.Range("K3").Value2 = 19.66432046836823
$ws.Range("L3").Value2 = 9.75040920902345
$ws.Range("O3").Value2 = 31.02074443741866
$ws.Range("C4").Value2 = 13.63440695200443
$ws.Range("D4").Value2 = 8.307245752697641
$ws.Range("E4").Value2 = 14.32426825826889
$ws.Range("F4").Value2 = 40.45138980958748
$ws.Range("G4").Value2 = 3.710357454841024
$ws.Range("J4").Value2 = 11.06839441863489
$ws.Range("K4").Value2 = 19.32448111909322
$ws.Range("L4").Value2 = 9.758115824532684
$ws.Range("O4").Value2 = 31.13359730611257
$ws.Range("C5").Value2 = 13.62443313087085
$ws.Range("D5").Value2 = 8.29508095480131
$ws.Range("E5").Value2 = 14.32416164181562
$ws.Range("F5").Value2 = 40.503431084345
$ws.Range("G5").Value2 = 3.71105980016774
$ws.Range("J5").Value2 = 11.07584654718379
$ws.Range("K5").Value2 = 19.18447524874149
$ws.Range("L5").Value2 = 9.761455976660468
$ws.Range("O5").Value2 = 31.18156419314714
$ws.Range("C6").Value2 = 13.622818597141
$ws.Range("D6").Value2 = 8.293070955201177
$ws.Range("E6").Value2 = 14.32417594985598
$ws.Range("F6").Value2 = 40.5122361665869
$ws.Range("G6").Value2 = 3.71117768581254
$ws.Range("J6").Value2 = 11.07710417554428
$ws.Range("K6").Value2 = 19.16114119390257
$ws.Range("L6").Value2 = 9.762022676909531
$ws.Range("O6").Value2 = 31.1896484119713
$ws.Range("C7").Value2 = 13.63426965638532
$ws.Range("D7").Value2 = 8.307081031524662
$ws.Range("E7").Value2 = 14.32426467530731
$ws.Range("F7").Value2 = 40.45208067801611
$ws.Range("G7").Value2 = 3.710366842375739
$ws.Range("J7").Value2 = 11.06849356611914
$ws.Range("K7").Value2 = 19.32259885125709
$ws.Range("L7").Value2 = 9.758160062074033
$ws.Range("O7").Value2 = 31.13423620008501
$ws.Range("C8").Value2 = 13.6907989025813
$ws.Range("D8").Value2 = 8.370296759764601
$ws.Range("E8").Value2 = 14.32929863225156
$ws.Range("F8").Value2 = 40.21048710769487
$ws.Range("G8").Value2 = 3.706971686816726
$ws.Range("J8").Value2 = 11.03343215763097
$ws.Range("K8").Value2 = 20.02125468013094
$ws.Range("L8").Value2 = 9.742866741976449
$ws.Range("O8").Value2 = 30.90709604144738
$ws.Range("C9").Value2 = 13.8226782393275
$ws.Range("D9").Value2 = 8.498731139873673
$ws.Range("E9").Value2 = 14.35574136882105
$ws.Range("F9").Value2 = 39.82313439128931
$ws.Range("G9").Value2 = 3.700966546785728
$ws.Range("J9").Value2 = 10.97520738617878
$ws.Range("K9").Value2 = 21.33028875073515
$ws.Range("L9").Value2 = 9.719143701334792
$ws.Range("O9").Value2 = 30.52447033322661
$ws.Range("C10").Value2 = 13.93185172392558
$ws.Range("D10").Value2 = 8.595323649259392
$ws.Range("E10").Value2 = 14.38515521969775
$ws.Range("F10").Value2 = 39.5918165190733
$ws.Range("G10").Value2 = 3.696948048866713
$ws.Range("J10").Value2 = 10.93883150893425
$ws.Range("K10").Value2 = 22.24581352327941
$ws.Range("L10").Value2 = 9.705509716281178
$ws.Range("O10").Value2 = 30.28193736340311
$ws.Range("C11").Value2 = 13.98405622388152
$ws.Range("D11").Value2 = 8.639643031232199
$ws.Range("E11").Value2 = 14.40067572352278
$ws.Range("F11").Value2 = 39.49827473104123
$ws.Range("G11").Value2 = 3.695204420259807
$ws.Range("J11").Value2 = 10.92366951031722
$ws.Range("K11").Value2 = 22.65061265166675
$ws.Range("L11").Value2 = 9.70012593279259
$ws.Range("O11").Value2 = 30.18004734337851
$ws.Range("C12").Value2 = 14.00417828042089
$ws.Range("D12").Value2 = 8.656470860748279
$ws.Range("E12").Value2 = 14.40685778876149
$ws.Range("F12").Value2 = 39.46454326261797
$ws.Range("G12").Value2 = 3.694556218798408
$ws.Range("J12").Value2 = 10.91812701026405
$ws.Range("K12").Value2 = 22.802098168058
$ws.Range("L12").Value2 = 9.69820445355535
$ws.Range("O12").Value2 = 30.14268369229568
$ws.Range("C13").Value2 = 13.99982911802702
$ws.Range("D13").Value2 = 8.652844834405432
$ws.Range("E13").Value2 = 14.4055128641039
$ws.Range("F13").Value2 = 39.47173257754324
$ws.Range("G13").Value2 = 3.694695284663287
$ws.Range("J13").Value2 = 10.91931183893567
$ws.Range("K13").Value2 = 22.76955509473996
$ws.Range("L13").Value2 = 9.698613071164838
$ws.Range("O13").Value2 = 30.15067626659748
$ws.Range("C14").Value2 = 13.98570466696502
$ws.Range("D14").Value2 = 8.641026618526205
$ws.Range("E14").Value2 = 14.40117823499639
$ws.Range("F14").Value2 = 39.4954656786088
$ws.Range("G14").Value2 = 3.695150850746803
$ws.Range("J14").Value2 = 10.92320953813524
$ws.Range("K14").Value2 = 22.66311226169915
$ws.Range("L14").Value2 = 9.699965504364416
$ws.Range("O14").Value2 = 30.17694892842285
$ws.Range("C15").Value2 = 13.97709868195129
$ws.Range("D15").Value2 = 8.633793204317614
$ws.Range("E15").Value2 = 14.39856275072049
$ws.Range("F15").Value2 = 39.5102233809785
$ws.Range("G15").Value2 = 3.695431468457012
$ws.Range("J15").Value2 = 10.92562290302735
$ws.Range("K15").Value2 = 22.59767465554922
$ws.Range("L15").Value2 = 9.700809164041418
$ws.Range("O15").Value2 = 30.19320074692003
$ws.Range("C16").Value2 = 13.92849015121183
$ws.Range("D16").Value2 = 8.592434092387926
$ws.Range("E16").Value2 = 14.38418373409472
$ws.Range("F16").Value2 = 39.59816573543603
$ws.Range("G16").Value2 = 3.697063692029765
$ws.Range("J16").Value2 = 10.93985024333993
$ws.Range("K16").Value2 = 22.21911310390175
$ws.Range("L16").Value2 = 9.705877989726185
$ws.Range("O16").Value2 = 30.28876647554971
$ws.Range("C17").Value2 = 13.89931280488355
$ws.Range("D17").Value2 = 8.567151972369546
$ws.Range("E17").Value2 = 14.3759087553029
$ws.Range("F17").Value2 = 39.65511599677351
$ws.Range("G17").Value2 = 3.698086580649607
$ws.Range("J17").Value2 = 10.94893297497287
$ws.Range("K17").Value2 = 21.98379501556229
$ws.Range("L17").Value2 = 9.709196847101238
$ws.Range("O17").Value2 = 30.34955895031322
$ws.Range("C18").Value2 = 13.88277069013176
$ws.Range("D18").Value2 = 8.552646605960559
$ws.Range("E18").Value2 = 14.37135082311251
$ws.Range("F18").Value2 = 39.68897160390625
$ws.Range("G18").Value2 = 3.698682867683541
$ws.Range("J18").Value2 = 10.95428755185344
$ws.Range("K18").Value2 = 21.84735375880152
$ws.Range("L18").Value2 = 9.711182823051857
$ws.Range("O18").Value2 = 30.38531914987602
$ws.Range("C19").Value2 = 13.87721137869135
$ws.Range("D19").Value2 = 8.547741849489594
$ws.Range("E19").Value2 = 14.36984229942907
$ws.Range("F19").Value2 = 39.70062302128877
$ws.Range("G19").Value2 = 3.698886127377109
$ws.Range("J19").Value2 = 10.9561229278343
$ws.Range("K19").Value2 = 21.80097342129509
$ws.Range("L19").Value2 = 9.711868488292385
$ws.Range("O19").Value2 = 30.39756311491619
$ws.Range("C20").Value2 = 13.90239403746365
$ws.Range("D20").Value2 = 8.569839614776805
$ws.Range("E20").Value2 = 14.37676879495152
$ws.Range("F20").Value2 = 39.64893969256984
$ws.Range("G20").Value2 = 3.697976870233227
$ws.Range("J20").Value2 = 10.94795260529328
$ws.Range("K20").Value2 = 22.00895903383505
$ws.Range("L20").Value2 = 9.708835577127415
$ws.Range("O20").Value2 = 30.34300526866667
$ws.Range("C21").Value2 = 13.98984386965572
$ws.Range("D21").Value2 = 8.644496768190052
$ws.Range("E21").Value2 = 14.40244317334163
$ws.Range("F21").Value2 = 39.48844872287693
$ws.Range("G21").Value2 = 3.695016712804577
$ws.Range("J21").Value2 = 10.92205928984806
$ws.Range("K21").Value2 = 22.694426967707
$ws.Range("L21").Value2 = 9.699565083556115
$ws.Range("O21").Value2 = 30.16919885255355
$ws.Range("C22").Value2 = 14.04905098762871
$ws.Range("D22").Value2 = 8.693547654862629
$ws.Range("E22").Value2 = 14.42099787691927
$ws.Range("F22").Value2 = 39.39341837903889
$ws.Range("G22").Value2 = 3.693152418895568
$ws.Range("J22").Value2 = 10.90629644608893
$ws.Range("K22").Value2 = 23.13186055570903
$ws.Range("L22").Value2 = 9.694189451534694
$ws.Range("O22").Value2 = 30.06271975926187
$ws.Range("C23").Value2 = 14.01726724994086
$ws.Range("D23").Value2 = 8.667347790229927
$ws.Range("E23").Value2 = 14.41093349752751
$ws.Range("F23").Value2 = 39.44323232757365
$ws.Range("G23").Value2 = 3.694141012779194
$ws.Range("J23").Value2 = 10.9146033125996
$ws.Range("K23").Value2 = 22.89939812750949
$ws.Range("L23").Value2 = 9.696996164382506
$ws.Range("O23").Value2 = 30.11889666202497
$ws.Range("C24").Value2 = 13.90100028764245
$ws.Range("D24").Value2 = 8.568624438560304
$ws.Range("E24").Value2 = 14.37637934950845
$ws.Range("F24").Value2 = 39.6517285298856
$ws.Range("G24").Value2 = 3.698026444718706
$ws.Range("J24").Value2 = 10.94839541676007
$ws.Range("K24").Value2 = 21.99758596532736
$ws.Range("L24").Value2 = 9.708998664534754
$ws.Range("O24").Value2 = 30.34596566603908
$ws.Range("C25").Value2 = 13.78480302128825
$ws.Range("D25").Value2 = 8.463558051262082
$ws.Range("E25").Value2 = 14.3468249272696
$ws.Range("F25").Value2 = 39.91861092699267
$ws.Range("G25").Value2 = 3.702521676984921
$ws.Range("J25").Value2 = 10.98983320408868
$ws.Range("K25").Value2 = 20.9836730909291
$ws.Range("L25").Value2 = 9.72489306215653
$ws.Range("O25").Value2 = 30.62122745328378
